# Update the "as_of_utc" timestamp column (AA) on both data sheets
# from "2025-11-04 08:10:46" to "2025-11-04 09:23:57" for rows 2-26.

$wb = $excel.ActiveWorkbook

$oldValue = "2025-11-04 08:10:46"
$newValue = "2025-11-04 09:23:57"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Range("AA$row")
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
